# Weekly update: add a new reporting date (2023-06-29, serial 45106) of
# price data for "Comercializadora del Agro de Limarí - Alcachofa", pushing
# the existing historical rows down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 3 new records right above the existing data block,
# which shifts rows 329:350 down to 332:353 (and extends the used range).
$ws.Rows("329:331").Insert()

# --- New row 329 ---
$ws.Range("A329").Value = 2
$ws.Range("B329").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C329").Value = 'Coquimbo'
$ws.Range("D329").Value = 45106
$ws.Range("E329").Value = 4
$ws.Range("F329").Value = 100112013
$ws.Range("G329").Value = 'Alcachofa'
$ws.Range("H329").Value = 'Argentina(o)'
$ws.Range("I329").Value = 'Primera'
$ws.Range("J329").Value = 1000
$ws.Range("K329").Value = 9000
$ws.Range("L329").Value = 10000
$ws.Range("M329").Value = 9500
$ws.Range("N329").Value = '$/caja 50 unidades'
$ws.Range("O329").Value = 'Provincia de Limarí'
$ws.Range("P329").Value = 190
$ws.Range("Q329").Value = 50
$ws.Range("R329").Value = 'Hortaliza'

# --- New row 330 ---
$ws.Range("A330").Value = 2
$ws.Range("B330").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C330").Value = 'Coquimbo'
$ws.Range("D330").Value = 45106
$ws.Range("E330").Value = 4
$ws.Range("F330").Value = 100112013
$ws.Range("G330").Value = 'Alcachofa'
$ws.Range("H330").Value = 'Española'
$ws.Range("I330").Value = 'Primera'
$ws.Range("J330").Value = 800
$ws.Range("K330").Value = 14000
$ws.Range("L330").Value = 15000
$ws.Range("M330").Value = 14500
$ws.Range("N330").Value = '$/caja 30 unidades'
$ws.Range("O330").Value = 'Provincia de Limarí'
$ws.Range("P330").Value = 483
$ws.Range("Q330").Value = 30
$ws.Range("R330").Value = 'Hortaliza'

# --- New row 331 ---
$ws.Range("A331").Value = 2
$ws.Range("B331").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C331").Value = 'Coquimbo'
$ws.Range("D331").Value = 45106
$ws.Range("E331").Value = 4
$ws.Range("F331").Value = 100112013
$ws.Range("G331").Value = 'Alcachofa'
$ws.Range("H331").Value = 'Madrigal'
$ws.Range("I331").Value = 'Primera'
$ws.Range("J331").Value = 360
$ws.Range("K331").Value = 12000
$ws.Range("L331").Value = 13000
$ws.Range("M331").Value = 12500
$ws.Range("N331").Value = '$/caja 40 unidades'
$ws.Range("O331").Value = 'Provincia de Limarí'
$ws.Range("P331").Value = 312
$ws.Range("Q331").Value = 40
$ws.Range("R331").Value = 'Hortaliza'
